# Added doors to board
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell -> new value (room-letter + door-direction suffix)
$doors = @{
    "E1"  = "OL"
    "H1"  = "OR"
    "O1"  = "SR"
    "V1"  = "KR"
    "A2"  = "AU"
    "Q3"  = "KL"
    "L6"  = "SD"
    "S6"  = "KD"
    "T6"  = "KD"
    "H8"  = "OR"
    "O10" = "DR"
    "A11" = "AD"
    "R14" = "CU"
    "S14" = "CU"
    "X14" = "LU"
    "E15" = "BU"
    "H15" = "HU"
    "N15" = "DD"
    "O18" = "CL"
    "S20" = "CR"
    "E22" = "BR"
    "L22" = "HR"
    "U22" = "LL"
}

foreach ($addr in $doors.Keys) {
    $ws.Range($addr).Value = $doors[$addr]
}

# Update the active selection on the sheet
$ws.Range("Q3").Select()
